# streamdetailedreport.xlsx — widen the two "Source" value columns and
# clarify their headers, per the commit:
#   "Added the description to header 'Value in SRC_DESC' instead of
#   'SRC_DESC' to detailed report."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("1st Source" values) needs to grow from ~14.76 to ~23.23
# characters wide, column E ("2nd Source" values) from ~15.52 to ~23.88,
# and column F (Difference) from ~14.754 to ~14.758 (essentially
# unchanged, rounding artifact of the original report generator).
# (Excel's COM layer stores column width using a 5-pixel cell-padding
# offset, so the ColumnWidth we assign is the target minus that offset;
# Excel itself re-adds it when it persists the <col> width to the sheet
# XML.)
$ws.Columns.Item(4).ColumnWidth = 22.401041666666668
$ws.Columns.Item(5).ColumnWidth = 23.041666666666668
$ws.Columns.Item(6).ColumnWidth = 13.924479166666666

# Update the header text for the two source-value columns to include the
# "Value in" prefix.
$ws.Range("D2").Value = "Value in 1st Source"
$ws.Range("E2").Value = "Value in 2nd Source"
